# Applies the "Updated symbol list" data refresh: refreshed Price/Volume(1h)
# figures for existing coins, plus a 6-row downward shift in the exchange-token
# block (rows 10-15) where WazirX/MandalaExchangeToken/BitrueCoin/BitMartToken/
# BitForexToken/One each moved up one row and "One" wrapped back in at row 15.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'254.84"
$ws.Range("E2").Value = "'3.69%"
$ws.Range("D3").Value = "'27.37"
$ws.Range("E3").Value = "'-9.16%"
$ws.Range("D4").Value = "'5.228"
$ws.Range("E4").Value = "'1.39%"
$ws.Range("D5").Value = "'0.05875"
$ws.Range("E5").Value = "'2.05%"
$ws.Range("E6").Value = "'0.77%"
$ws.Range("D7").Value = "'3.220"
$ws.Range("E7").Value = "'-2.21%"
$ws.Range("D8").Value = "'0.8650"
$ws.Range("E8").Value = "'1.73%"
$ws.Range("D9").Value = "'0.9665"
$ws.Range("E9").Value = "'12.71%"
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D10").Value = "'0.1412"
$ws.Range("E10").Value = "'2.05%"
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D11").Value = "'0.07162"
$ws.Range("E11").Value = "'1.16%"
$ws.Range("B12").Value = "BitrueCoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D12").Value = "'0.03185"
$ws.Range("E12").Value = "'-1.44%"
$ws.Range("B13").Value = "BitMartToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D13").Value = "'0.09233"
$ws.Range("E13").Value = "'-1.39%"
$ws.Range("B14").Value = "BitForexToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D14").Value = "'0.001554"
$ws.Range("E14").Value = "'1.19%"
$ws.Range("B15").Value = "One"
$ws.Range("C15").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D15").Value = "'0.0006100"
$ws.Range("E15").Value = "'2.17%"
$ws.Range("D16").Value = "'0.005805"
$ws.Range("E16").Value = "'-3.51%"
$ws.Range("D17").Value = "'3.501"
$ws.Range("E17").Value = "'-1.34%"
$ws.Range("E18").Value = "'1.85%"
$ws.Range("E20").Value = "'1.95%"
$ws.Range("D21").Value = "'0.1308"
$ws.Range("E21").Value = "'-1.17%"
$ws.Range("D22").Value = "'3.526"
$ws.Range("E22").Value = "'1.01%"
$ws.Range("D23").Value = "'0.04148"
$ws.Range("E23").Value = "'0.35%"
$ws.Range("E24").Value = "'-2.13%"
$ws.Range("D25").Value = "'0.001223"
$ws.Range("E25").Value = "'-0.03%"
$ws.Range("D26").Value = "'0.004804"
$ws.Range("E26").Value = "'15.56%"
$ws.Range("E27").Value = "'-0.02%"
$ws.Range("E28").Value = "'1.18%"
$ws.Range("D40").Value = "'0.03813"
$ws.Range("E40").Value = "'1.73%"
$ws.Range("D41").Value = "'0.005669"
$ws.Range("E41").Value = "'-0.20%"
$ws.Range("E42").Value = "'3.07%"
$ws.Range("D43").Value = "'0.002343"
$ws.Range("E43").Value = "'11.57%"
$ws.Range("D44").Value = "'0.01064"
$ws.Range("E44").Value = "'-0.56%"
$ws.Range("D45").Value = "'0.00005241"
$ws.Range("E45").Value = "'-4.43%"
$ws.Range("E46").Value = "'0.00%"
$ws.Range("D47").Value = "'0.10000"
$ws.Range("E47").Value = "'40.84%"
$ws.Range("D48").Value = "'0.002132"
$ws.Range("E48").Value = "'-14.03%"
$ws.Range("E49").Value = "'0.00%"
$ws.Range("E50").Value = "'0.00%"
